$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Last status check on" timestamp in F1
$ws.Range("F1").Value = "Last status check on: 25.02.2022 13:00"

# Row 7 (MOL Olomoucká): price moved up, old price shifted into "Old Cena",
# delta now shown as text with explicit sign, and date shown as a text timestamp
$ws.Range("B7").Value = 38.5
$ws.Range("C7").Value = 38.29

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "+0.21"
$ws.Range("D7").ClearFormats()

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "2022-02-25 13:00:19"
$ws.Range("E7").ClearFormats()
